# TC26_Trials_Filter_Diagnosis-UterineCancer.xlsx
# - add the Neo4j MATCH query text to cell A2 on the "startup" sheet
# - grow row 2 to fit the wrapped query text
# - update the sheet selection to B2:B5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$query = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Uterine cancer, NOS'] RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(t.clinical_trial_designation ,'')as ``Trial Code`` , coalesce(a.arm_id,'') As ``Arm`` , coalesce(a.arm_drug,'') As ``Arm Treatment`` , coalesce(c.disease,'') As Diagnosis , coalesce(c.gender,'') As Gender , coalesce(c.race,'') As Race , coalesce(c.ethnicity,'') As Ethnicity"

$ws.Range("A2").Value = $query

# Row 2 grows tall enough to show the wrapped query text
$ws.Rows.Item(2).RowHeight = 87

# Select B2:B5 as the active range on the sheet
$ws.Range("B2:B5").Select()
